# Daily update: Feb 6, 2026 - 166 new games, 20091 season games
# Swap the C:K (firstName..steals) contents of the listed row pairs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(52, 53),
    @(297, 298),
    @(313, 314),
    @(332, 333),
    @(363, 364),
    @(431, 432)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("C$($r1):K$($r1)")
    $range2 = $ws.Range("C$($r2):K$($r2)")

    $vals1 = $range1.Value()
    $vals2 = $range2.Value()

    $range1.Value = $vals2
    $range2.Value = $vals1
}
